# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (quarterly fund-holding detail, same
#    layout as the other quarter sheets) right before the "总计" (total)
#    summary sheet.
# 2. Add a new top row to the "总计" sheet summarizing the 2022-Q1 quarter
#    (holding count + market value), pushing the previous rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Locate the existing sheets we need as references / anchors.
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right before "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Re-resolve "总计" AFTER the insert: adding a sheet in front of it shifts
# its position, and a Worksheet handle captured beforehand tracks the old
# *position* rather than the sheet's identity, so it must be re-fetched.
$totalSheet = $wb.Worksheets.Item("总计")

# Match the page-margin conventions used by the other quarter sheets.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Copy the header-row / index-column formatting (bold, centered, bordered)
# from the 2021-Q4 sheet so the new sheet matches the established style.
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Helper: write a text value without Excel's "looks-like-a-number" auto
# conversion (leading apostrophe forces text), then drop back to the
# default "Normal" style so no stray per-cell style (e.g. quote-prefix)
# is left behind.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2: 007012 / 湘财长顺混合A
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "007012"
Set-TextValue $newSheet.Range("C2") "湘财长顺混合A"
Set-TextValue $newSheet.Range("D2") "4.70"
Set-TextValue $newSheet.Range("E2") "94.08"
Set-TextValue $newSheet.Range("F2") "6.66"
Set-TextValue $newSheet.Range("G2") "0.3130"
$newSheet.Range("H2").Value = 4

# Row 3: 008128 / 湘财长源股票A
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "008128"
Set-TextValue $newSheet.Range("C3") "湘财长源股票A"
Set-TextValue $newSheet.Range("D3") "2.74"
Set-TextValue $newSheet.Range("E3") "94.29"
Set-TextValue $newSheet.Range("F3") "6.63"
Set-TextValue $newSheet.Range("G3") "0.1817"
$newSheet.Range("H3").Value = 5

# Row 4: 007013 / 湘财长顺混合C
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "007013"
Set-TextValue $newSheet.Range("C4") "湘财长顺混合C"
Set-TextValue $newSheet.Range("D4") "2.47"
Set-TextValue $newSheet.Range("E4") "94.08"
Set-TextValue $newSheet.Range("F4") "6.66"
Set-TextValue $newSheet.Range("G4") "0.1645"
$newSheet.Range("H4").Value = 4

# Row 5: 008129 / 湘财长源股票C
$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet.Range("B5") "008129"
Set-TextValue $newSheet.Range("C5") "湘财长源股票C"
Set-TextValue $newSheet.Range("D5") "1.05"
Set-TextValue $newSheet.Range("E5") "94.29"
Set-TextValue $newSheet.Range("F5") "6.63"
Set-TextValue $newSheet.Range("G5") "0.0696"
$newSheet.Range("H5").Value = 5

# Row 6: 005290 / 诺德新盛灵活配置混合A
$newSheet.Range("A6").Value = 4
Set-TextValue $newSheet.Range("B6") "005290"
Set-TextValue $newSheet.Range("C6") "诺德新盛灵活配置混合A"
Set-TextValue $newSheet.Range("D6") "0.05"
Set-TextValue $newSheet.Range("E6") "91.12"
Set-TextValue $newSheet.Range("F6") "1.48"
Set-TextValue $newSheet.Range("G6") "0.0007"
$newSheet.Range("H6").Value = 10

# Row 7: 009710 / 诺德新盛灵活配置混合C
$newSheet.Range("A7").Value = 5
Set-TextValue $newSheet.Range("B7") "009710"
Set-TextValue $newSheet.Range("C7") "诺德新盛灵活配置混合C"
Set-TextValue $newSheet.Range("D7") "0.03"
Set-TextValue $newSheet.Range("E7") "91.12"
Set-TextValue $newSheet.Range("F7") "1.48"
Set-TextValue $newSheet.Range("G7") "0.0004"
$newSheet.Range("H7").Value = 10

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing quarter rows down by one
#    and insert the new 2022-Q1 summary row at the top (row 2).
# ---------------------------------------------------------------------

# Give the new bottom row (row 6, previously unused) the same index-column
# style ("s=2") the other data rows already carry.
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q1"
$totalSheet.Range("C6").Value = 1
$totalSheet.Range("D6").Value = 0.07000000000000001

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.07000000000000001

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.01

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 5
$totalSheet.Range("D3").Value = 0.95

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.73

# ---------------------------------------------------------------------
# Restore the original active sheet/tab (adding a sheet shifts focus).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
